$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 39, shifting the existing
# rows 39-48 down to 41-50 (they keep their original values).
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()

# New row 39: Damasco, Primera, Región Metropolitana, $/caja 16 kilos
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44900
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100103
$ws.Cells.Item(39, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39, 9).Value = 100103003
$ws.Cells.Item(39, 10).Value = "Damasco"
$ws.Cells.Item(39, 11).Value = "Castle Brite"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 200
$ws.Cells.Item(39, 14).Value = 23000
$ws.Cells.Item(39, 15).Value = 24000
$ws.Cells.Item(39, 16).Value = 23500
$ws.Cells.Item(39, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(39, 18).Value = "Región Metropolitana"
$ws.Cells.Item(39, 19).Value = 1469
$ws.Cells.Item(39, 20).Value = 16

# New row 40: Damasco, Segunda, Región Metropolitana, $/caja 16 kilos
$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 44900
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100103
$ws.Cells.Item(40, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(40, 9).Value = 100103003
$ws.Cells.Item(40, 10).Value = "Damasco"
$ws.Cells.Item(40, 11).Value = "Castle Brite"
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value = 100
$ws.Cells.Item(40, 14).Value = 19000
$ws.Cells.Item(40, 15).Value = 19000
$ws.Cells.Item(40, 16).Value = 19000
$ws.Cells.Item(40, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(40, 18).Value = "Región Metropolitana"
$ws.Cells.Item(40, 19).Value = 1188
$ws.Cells.Item(40, 20).Value = 16
